$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Fix capitalization of "de" -> "De" in a couple of data cells
$ws.Range("A2").Value = "Ciudad De México"
$ws.Range("B8").Value = "Izúcar De Matamoros"

# Remove the trailing footer/metadata rows (14-18)
$ws.Range("A14:D18").EntireRow.Delete()
